$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 10104
$ws.Range("D2").Value = 9994

$ws.Range("C3").Value = 10104
$ws.Range("D3").Value = 9994

$ws.Range("C4").Value = 10104
$ws.Range("D4").Value = 9991

$ws.Range("C5").Value = 10104
$ws.Range("D5").Value = 9994

$ws.Range("C6").Value = 4966
$ws.Range("D6").Value = 9971

$ws.Range("C7").Value = 10104
$ws.Range("D7").Value = 9994

$ws.Range("C8").Value = 10104
$ws.Range("D8").Value = 5064

$ws.Range("C9").Value = 5911
$ws.Range("D9").Value = 6975

$ws.Range("C10").Value = 5911
$ws.Range("D10").Value = 4440
